$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder merge cells: move A5:A7 ahead of A1:A4 and A8:A10
# by unmerging and re-merging A1:A4 then A8:A10 (A5:A7 stays put,
# which bumps it to the front of the merge list).
$ws.Range("A1:A4").UnMerge()
$ws.Range("A8:A10").UnMerge()
$ws.Range("A1:A4").Merge()
$ws.Range("A8:A10").Merge()

# Update the refreshed MAA percentage figures / timestamp
$ws.Range('D2').Value = 'maa://24702 (94.12), maa://25390 (95.78), maa://36681 (88.24)'
$ws.Range('AB2').Value = 'maa://21246 (91.26), maa://36684 (98.7), ***maa://22731 (6.67)'
$ws.Range('L3').Value = '*maa://22880 (69.57), maa://20276 (83.33), *maa://22749 (66.67)'
$ws.Range('X3').Value = 'maa://27396 (85.42), maa://27484 (95.79), maa://27480 (82.35)'
$ws.Range('AB3').Value = 'maa://24390 (96.23)'
$ws.Range('D4').Value = 'maa://24632 (93.48), **maa://24303 (36.36), maa://22499 (85.71), maa://22746 (100.0)'
$ws.Range('T4').Value = 'maa://32509 (97.7), maa://22754 (91.67), maa://27295 (81.82), *maa://21746 (55.81), *maa://31008 (78.05)'
$ws.Range('D5').Value = 'maa://21245 (82.14), maa://22744 (83.33)'
$ws.Range('D6').Value = 'maa://42407 (94.44)'
$ws.Range('P6').Value = 'maa://31836 (88.89), maa://30381 (91.67)'
$ws.Range('D7').Value = 'maa://21955 (93.55)'
$ws.Range('P7').Value = 'maa://22750 (94.74)'
$ws.Range('X7').Value = 'maa://22399 (94.78), *maa://22758 (70.37)'
$ws.Range('AF7').Value = '*maa://26191 (68.92), *maa://36671 (73.33), *maa://42530 (75.0)'
$ws.Range('A8').Value = '更新日期：2024.11.02 22:47:19'
$ws.Range('X9').Value = 'maa://26223 (97.12)'
$ws.Range('AB9').Value = 'maa://28711 (88.51), ***maa://22740 (5.88), **maa://27377 (46.15), ***maa://25174 (20.0), **maa://39938 (41.18), maa://40166 (100.0)'
$ws.Range('AF9').Value = 'maa://26206 (90.22), **maa://22865 (47.92)'
$ws.Range('X12').Value = 'maa://22753 (91.5), *maa://21485 (77.1), maa://37962 (83.33)'
$ws.Range('AB12').Value = 'maa://23669 (95.51), maa://36677 (92.68), maa://39872 (86.67)'
$ws.Range('AF12').Value = '*maa://28932 (78.15), *maa://20106 (63.64), *maa://22769 (64.29)'
$ws.Range('D13').Value = 'maa://24999 (91.45), maa://36673 (92.06), maa://25001 (85.51)'
$ws.Range('H13').Value = '*maa://21248 (74.53), **maa://22728 (46.51)'
$ws.Range('X13').Value = '*maa://34957 (77.78), *maa://22768 (51.61)'
$ws.Range('AF13').Value = '**maa://22737 (30.6), maa://39883 (91.43), *maa://39885 (70.0)'
$ws.Range('L14').Value = 'maa://26245 (96.12), maa://21288 (96.21), maa://36682 (100.0), maa://39841 (93.75)'
$ws.Range('D15').Value = '*maa://22743 (77.09), maa://22734 (83.48), *maa://30808 (64.41), ***maa://36048 (12.12)'
$ws.Range('H15').Value = 'maa://24304 (88.3), maa://21478 (91.18)'
$ws.Range('P15').Value = 'maa://24762 (89.73), *maa://22727 (70.0)'
$ws.Range('AF15').Value = 'maa://21364 (80.68), *maa://22766 (72.12), *maa://36666 (77.94)'
$ws.Range('T16').Value = 'maa://22729 (95.21), *maa://28648 (69.64), maa://36674 (81.25)'
$ws.Range('X16').Value = 'maa://28501 (97.67), maa://28051 (95.83)'
$ws.Range('T17').Value = '***maa://42324 (28.57)'
$ws.Range('T19').Value = 'maa://24386 (98.85)'
$ws.Range('AB19').Value = '*maa://30709 (61.56), *maa://36668 (52.17)'
$ws.Range('H20').Value = 'maa://22864 (88.89)'
$ws.Range('L20').Value = 'maa://41331 (82.46)'
$ws.Range('H21').Value = 'maa://24372 (96.51)'
$ws.Range('AF21').Value = 'maa://22524 (94.02), *maa://22432 (75.44)'
$ws.Range('T22').Value = 'maa://38495 (88.89)'
$ws.Range('X22').Value = 'maa://21282 (98.86), *maa://37649 (71.43)'
$ws.Range('L23').Value = 'maa://39756 (93.1), maa://39875 (94.34)'
$ws.Range('P23').Value = 'maa://30587 (91.33), *maa://29748 (75.2), ***maa://29785 (15.15), *maa://37566 (76.19)'
$ws.Range('AF24').Value = 'maa://22523 (85.34), *maa://36672 (77.78), maa://29910 (94.12), **maa://21440 (34.55)'
$ws.Range('D25').Value = 'maa://29753 (94.89)'
$ws.Range('H26').Value = 'maa://24913 (91.67)'
$ws.Range('AB26').Value = 'maa://42235 (87.18)'
$ws.Range('D28').Value = 'maa://24465 (90.48), maa://25725 (82.72)'
$ws.Range('X28').Value = 'maa://39929 (87.55), ***maa://39723 (14.29), maa://41749 (83.33)'
$ws.Range('AF28').Value = 'maa://36660 (92.7), *maa://36701 (62.96)'
$ws.Range('L29').Value = 'maa://28432 (93.6), *maa://28440 (72.84), maa://31400 (100.0), *maa://28650 (66.67)'
$ws.Range('T30').Value = '*maa://32940 (66.67), maa://24388 (94.12)'
$ws.Range('L31').Value = 'maa://35926 (93.36), *maa://36258 (79.75)'
$ws.Range('T32').Value = 'maa://41108 (88.64), maa://41238 (95.24), maa://42859 (92.86)'
$ws.Range('L35').Value = 'maa://41296 (98.59)'
$ws.Range('AF38').Value = 'maa://36697 (84.72)'
$ws.Range('P40').Value = 'maa://23278 (95.92), maa://21386 (95.68), maa://36664 (91.11)'
$ws.Range('H43').Value = 'maa://22525 (92.8), maa://21284 (82.93)'
$ws.Range('H44').Value = 'maa://29768 (97.65), maa://27728 (96.0)'
$ws.Range('H46').Value = 'maa://35931 (92.15)'
$ws.Range('H47').Value = 'maa://27410 (95.86), maa://29661 (97.73), maa://28038 (84.62)'
$ws.Range('H53').Value = 'maa://32534 (93.09), **maa://32434 (34.78)'
$ws.Range('H60').Value = '**maa://40438 (36.36)'
